$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the learning parameter values that were missing for the
# "input" (row 15) and "output" (row 16) dims.
$ws.Range("D15").Value = 1
$ws.Range("D16").Value = 0

# Update the selected cell to reflect where the user left off editing.
$ws.Range("D17").Select()
